$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "walkingToRunning"
$ws.Range("C22").Value = 20.3416576385498
$ws.Range("D22").Value = -37.60380554199219
$ws.Range("E22").Value = 9.662521362304688
$ws.Range("F22").Value = 2.631648543253113
$ws.Range("G22").Value = -5.290064356919663
$ws.Range("H22").Value = -1.781815700103778

$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "walkingToRunning"
$ws.Range("C23").Value = -8.722187995910645
$ws.Range("D23").Value = 8.096003532409668
$ws.Range("E23").Value = -8.941717147827148
$ws.Range("F23").Value = 0.2861763049412409
$ws.Range("G23").Value = -0.7290360459013701
$ws.Range("H23").Value = -0.5830237417551998

$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "walkingToRunning"
$ws.Range("C24").Value = -72.57207489013672
$ws.Range("D24").Value = -62.40628051757812
$ws.Range("E24").Value = 63.55230712890625
$ws.Range("F24").Value = -4.221912179034594
$ws.Range("G24").Value = -0.03346766901844722
$ws.Range("H24").Value = 9.115482657631061

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "walkingToRunning"
$ws.Range("C25").Value = -86.77053833007812
$ws.Range("D25").Value = -16.04696655273438
$ws.Range("E25").Value = -37.80954742431641
$ws.Range("F25").Value = -1.58611540711694
$ws.Range("G25").Value = 0.8197227957620261
$ws.Range("H25").Value = 4.598585906056506

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "walkingToRunning"
$ws.Range("C26").Value = -12.0439281463623
$ws.Range("D26").Value = -10.04275798797607
$ws.Range("E26").Value = -17.28581237792969
$ws.Range("F26").Value = 2.200832122323145
$ws.Range("G26").Value = -1.895753322998213
$ws.Range("H26").Value = -4.244981093213744

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "walkingToRunning"
$ws.Range("C27").Value = -29.28538513183594
$ws.Range("D27").Value = -21.46180725097656
$ws.Range("E27").Value = 6.313263416290283
$ws.Range("F27").Value = 2.756207520217482
$ws.Range("G27").Value = -7.966393658191489
$ws.Range("H27").Value = -2.264977694935836

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "walkingToRunning"
$ws.Range("C28").Value = 3.212208271026612
$ws.Range("D28").Value = 6.285176277160645
$ws.Range("E28").Value = 0.2726368904113769
$ws.Range("F28").Value = 2.778578069857753
$ws.Range("G28").Value = -3.013215610746667
$ws.Range("H28").Value = -2.033113320094297

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "walkingToRunning"
$ws.Range("C29").Value = -6.000577926635742
$ws.Range("D29").Value = 12.46978282928467
$ws.Range("E29").Value = 1.907239437103272
$ws.Range("F29").Value = -1.761591400025211
$ws.Range("G29").Value = 3.015217549539029
$ws.Range("H29").Value = 1.020017661111206

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "walkingToRunning"
$ws.Range("C30").Value = 19.4264965057373
$ws.Range("D30").Value = -45.61504745483398
$ws.Range("E30").Value = 8.111544609069824
$ws.Range("F30").Value = -4.130074517575332
$ws.Range("G30").Value = -0.597962478681886
$ws.Range("H30").Value = -1.069513860465531

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "walkingToRunning"
$ws.Range("C31").Value = -43.24179840087891
$ws.Range("D31").Value = 13.53559017181396
$ws.Range("E31").Value = -9.444896697998049
$ws.Range("F31").Value = -2.031331133980304
$ws.Range("G31").Value = 4.770669049610229
$ws.Range("H31").Value = -3.800996664631577

